$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (header in G1 is "K", replacing old "Strike#").
# Update rows 2-10 in column G per regenerated save_data / K calc.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 1
